$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data rows (rows 2-7), the new data will fully replace them.
$ws.Range("A2:T7").ClearContents()

# Full replacement data for rows 2-10: every combination of sending/target
# cluster across ECs, FAPs and MuSCs for the Ptn -> Plxnb2 ligand-receptor
# pair, recomputed with the updated TPM values.
$data = @(
    @("ECs","Ptn","Plxnb2","ECs",1,0.3333333333333333,0.06520933333333333,0.195628,0.007116089623850998,0.007116089623850999,3,1,2.325008666666667,6.975026,0.05445297772988467,0.05445297772988466,0.1516122651475555,1.364510386328,0.0003874922698114218,0.0003874922698114217),
    @("ECs","Ptn","Plxnb2","FAPs",1,0.3333333333333333,0.06520933333333333,0.195628,0.007116089623850998,0.007116089623850999,3,1,19.33828433333333,58.014853,0.4529132218878514,0.4529132218878514,1.261036629187111,11.349329662684,0.003222971078781064,0.003222971078781064),
    @("ECs","Ptn","Plxnb2","MuSCs",1,0.3333333333333333,0.06520933333333333,0.195628,0.007116089623850998,0.007116089623850999,3,1,21.03425566666667,63.102767,0.492633800382264,0.492633800382264,1.371629789186222,12.344668102676,0.003505626275258513,0.003505626275258513),
    @("FAPs","Ptn","Plxnb2","ECs",3,1,3.270036666666666,9.81011,0.3568488252184601,0.3568488252184602,3,1,2.325008666666667,6.975026,0.05445297772988467,0.05445297772988466,7.602863590317778,68.42577231285999,0.01943148113255631,0.01943148113255631),
    @("FAPs","Ptn","Plxnb2","FAPs",3,1,3.270036666666666,9.81011,0.3568488252184601,0.3568488252184602,3,1,19.33828433333333,58.014853,0.4529132218878514,0.4529132218878514,63.23689884042555,569.13208956383,0.1616215511565875,0.1616215511565876),
    @("FAPs","Ptn","Plxnb2","MuSCs",3,1,3.270036666666666,9.81011,0.3568488252184601,0.3568488252184602,3,1,21.03425566666667,63.102767,0.492633800382264,0.492633800382264,68.78278728604111,619.04508557437,0.1757957929293163,0.1757957929293163),
    @("MuSCs","Ptn","Plxnb2","ECs",3,1,5.828401,17.485203,0.6360350851576888,0.6360350851576889,3,1,2.325008666666667,6.975026,0.05445297772988467,0.05445297772988466,13.55108283780867,121.959745540278,0.03463400432751693,0.03463400432751693),
    @("MuSCs","Ptn","Plxnb2","FAPs",3,1,5.828401,17.485203,0.6360350851576888,0.6360350851576889,3,1,19.33828433333333,58.014853,0.4529132218878514,0.4529132218878514,112.7112757466844,1014.401481720159,0.2880686996524828,0.2880686996524828),
    @("MuSCs","Ptn","Plxnb2","MuSCs",3,1,5.828401,17.485203,0.6360350851576888,0.6360350851576889,3,1,21.03425566666667,63.102767,0.492633800382264,0.492633800382264,122.5960767618557,1103.364690856701,0.3133323811776891,0.3133323811776892)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}
